# Apply the data updates + header/title restyle described by the commit.
#
# 1) Training Dashboard sheet: "PERIOD TO EXPIRE" (H3) and "LAST UPDATE" (I3)
#    values are refreshed to reflect a later check/run.
# 2) The title row and the column-header row (on every sheet) are restyled so
#    the bold header text is white (instead of plain black) and uses the
#    workbook's normal font size -- this makes the bold white text read
#    clearly against the dark-blue header fill used across the dashboards.

$wb = $excel.ActiveWorkbook

$white = 16777215  # RGB(255,255,255)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $dim = $ws.UsedRange
    $lastCol = $dim.Columns.Count

    # Title cell (row 1) -- bold, normal size, white text.
    $titleRange = $ws.Range("A1")
    $titleRange.Font.Bold = $true
    $titleRange.Font.Size = 11
    $titleRange.Font.Color = $white

    # Column header row (row 2) -- bold, white text.
    $headerRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $lastCol))
    $headerRange.Font.Bold = $true
    $headerRange.Font.Color = $white
}

# Training Dashboard data refresh.
$trainingWs = $wb.Worksheets.Item("Training Dashboard")
$trainingWs.Range("H3").Value = -55

$trainingWs.Range("I3").NumberFormat = "@"
$trainingWs.Range("I3").Value = "16-Sep-2025"
